$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new players in rows 28 and 29 -------------------------------
# Row 28: Kamal Bafounta, 23 y/o, born 1/8/2002 (serial 37264), poste MC
$ws.Range("A28").Value = "Kamal Bafounta"
$ws.Range("A28").HorizontalAlignment = -4108
$ws.Range("A28").VerticalAlignment = -4108
$ws.Range("B28").Value = 23
$ws.Range("C28").Value = 37264
$ws.Range("D28").Value = "MC"

# Row 29: Malik Boussaïd, 29 y/o, born 9/3/1996 (serial 35311), poste DD
$ws.Range("A29").Value = "Malik Boussaïd"
$ws.Range("A29").HorizontalAlignment = -4108
$ws.Range("A29").VerticalAlignment = -4108
$ws.Range("B29").Value = 29
$ws.Range("C29").Value = 35311
$ws.Range("D29").Value = "DD"

# --- Extend the conditional formatting range to cover the new rows -------
$fcs = $ws.Range("A22:A27").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A22:A29"))
}

# --- Update the sheet's last active selection -----------------------------
$ws.Range("C33").Select() | Out-Null
